# ZBP_11_obavy_epidemie.xlsx update:
#  - add a new date column "16. 3. 2021" (between the existing last date
#    column and the "total" summary column) on both the "data" and
#    "pocetR" worksheets, with its per-row values
#  - bump the "aktualizace 9. 3. 2021" -> "aktualizace 23. 3. 2021" in the
#    two summary footer cells

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "data": dates live in columns D..Z (col 4..26); the "total" /
# category columns A:C follow immediately after. Insert the new date
# column right before the old column AA (i.e. right after Z), which
# becomes the new column AA.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("data")

$wsData.Columns.Item(27).Insert()

$wsData.Cells.Item(1, 27).Value = "16. 3. 2021"

$dataValues = @(
    0.16, 0.35, 0.49, 0.24, 0.41, 0.35, 0.17, 0.36, 0.47, 0.11,
    0.31, 0.58, 0.17, 0.33, 0.5,  0.14, 0.31, 0.55, 0.18, 0.42,
    0.4,  0.23, 0.38, 0.39, 0.11, 0.32, 0.57, 0.16, 0.36, 0.48,
    0.17, 0.31, 0.52, 0.18, 0.36, 0.46, 0.14, 0.38, 0.48, 0.05,
    0.16, 0.79, 0.14, 0.35, 0.51, 0.15, 0.38, 0.47, 0.12, 0.41,
    0.47, 0.22, 0.37, 0.41, 0.2,  0.34, 0.46, 0.31, 0.39, 0.3,
    0.23, 0.4,  0.37, 0.17, 0.35, 0.48, 0.16, 0.43, 0.41, 0.12,
    0.31, 0.57, 0.07000000000000001, 0.28, 0.65
)

for ($i = 0; $i -lt $dataValues.Length; $i++) {
    $wsData.Cells.Item($i + 2, 27).Value = $dataValues[$i]
}

# Footer row: bump the date in the summary sentence.
$wsData.Cells.Item(77, 1).Value = "Život během pandemie, Obavy z epidemie, % respondentů celkově a ve skupinách, aktualizace 23. 3. 2021"

# ---------------------------------------------------------------------
# Sheet "pocetR": dates live in columns C..Y (col 3..25); insert the new
# date column right before the old column Z, which becomes the new
# column Z.
# ---------------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("pocetR")

$wsCounts.Columns.Item(26).Insert()

$wsCounts.Cells.Item(1, 26).Value = "16. 3. 2021"

$countValues = @(
    2101, 500, 776, 825, 676, 744, 681, 1015, 1086, 1104,
    478,  242, 277, 46,  154, 100, 22,  303,  563,  256,
    391,  368, 244, 385, 457
)

for ($i = 0; $i -lt $countValues.Length; $i++) {
    $wsCounts.Cells.Item($i + 2, 26).Value = $countValues[$i]
}

# Footer row has an (empty-text) cell in every trailing column; match that
# for the freshly-inserted column too by copying the same "blank" cell
# from its left neighbour.
$wsCounts.Cells.Item(27, 25).Copy($wsCounts.Cells.Item(27, 26))

# Footer row: bump the date in the summary sentence.
$wsCounts.Cells.Item(27, 1).Value = "Život během pandemie, Obavy z epidemie, velikost dotázaného souboru celkově a ve skupinách, aktualizace 23. 3. 2021"
